$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 28682.666
$ws.Range("I34").Value = 11499.5
$ws.Range("J34").Value = 63049
$ws.Range("K34").Value = 11499.5
$ws.Range("L34").Value = 63049
$ws.Range("M34").Value = -11296.5
$ws.Range("N34").Value = -63455
$ws.Range("H36").Value = 28682.666
$ws.Range("I36").Value = 11499.5
$ws.Range("J36").Value = 63049
$ws.Range("K36").Value = 11499.5
$ws.Range("L36").Value = 63049
$ws.Range("M36").Value = -10784.5
$ws.Range("N36").Value = -64479
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").ClearContents()
$ws.Range("N49").Value = 0
$ws.Range("H70").Value = 65000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 65000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 195000
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -195540
$ws.Range("H73").Value = 65000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 65000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 195000
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -196872
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("H112").Value = 1391.2559
$ws.Range("I112").Value = 999
$ws.Range("J112").Value = 1400.5952
$ws.Range("K112").Value = 2997
$ws.Range("L112").Value = 4201.7856
$ws.Range("M112").Value = -1889
$ws.Range("N112").Value = -6417.7856
$ws.Range("H113").Value = 26391.875
$ws.Range("I113").Value = 41301
$ws.Range("J113").Value = 1543.3334
$ws.Range("K113").Value = 41301
$ws.Range("L113").Value = 1543.3334
$ws.Range("M113").Value = -38047
$ws.Range("N113").Value = -8051.3334
$ws.Range("H132").Value = 1081.8572
$ws.Range("I132").Value = 971.8182
$ws.Range("J132").Value = 1485.3334
$ws.Range("K132").Value = 2915.4546
$ws.Range("L132").Value = 4456.0002
$ws.Range("M132").Value = -385.4546
$ws.Range("N132").Value = -9516.0002
$ws.Range("H137").Value = 2375
$ws.Range("I137").Value = 1300
$ws.Range("J137").Value = 3450
$ws.Range("K137").Value = 3900
$ws.Range("L137").Value = 10350
$ws.Range("M137").Value = -1350
$ws.Range("N137").Value = -15450

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1303.7142
$ws.Range("I45").Value = 923.8333
$ws.Range("J45").Value = 1810.2222
$ws.Range("K45").Value = 923.8333
$ws.Range("L45").Value = 1810.2222
$ws.Range("M45").Value = -546.8333
$ws.Range("N45").Value = -2564.2222
$ws.Range("H74").Value = 1649.3478
$ws.Range("I74").Value = 485
$ws.Range("J74").Value = 4310.7144
$ws.Range("K74").Value = 485
$ws.Range("L74").Value = 4310.7144
$ws.Range("M74").Value = 389
$ws.Range("N74").Value = -6058.7144
$ws.Range("H77").Value = 1649.3478
$ws.Range("I77").Value = 485
$ws.Range("J77").Value = 4310.7144
$ws.Range("K77").Value = 2425
$ws.Range("L77").Value = 21553.572
$ws.Range("M77").Value = 1943
$ws.Range("N77").Value = -30289.572
$ws.Range("H102").Value = 1416
$ws.Range("I102").Value = 1177.2
$ws.Range("J102").Value = 1615
$ws.Range("K102").Value = 1177.2
$ws.Range("L102").Value = 1615
$ws.Range("M102").Value = 444.8
$ws.Range("N102").Value = -4859
$ws.Range("H132").Value = 2226.9048
$ws.Range("I132").Value = 1876
$ws.Range("J132").Value = 4332.3335
$ws.Range("K132").Value = 5628
$ws.Range("L132").Value = 12997.0005
$ws.Range("M132").Value = -3098
$ws.Range("N132").Value = -18057.0005
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").ClearContents()
$ws.Range("N138").Value = 0

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 103063.3
$ws.Range("I86").Value = 3479.0667
$ws.Range("J86").Value = 401816
$ws.Range("K86").Value = 3479.0667
$ws.Range("L86").Value = 401816
$ws.Range("M86").Value = -2356.0667
$ws.Range("N86").Value = -404062
$ws.Range("H89").Value = 103063.3
$ws.Range("I89").Value = 3479.0667
$ws.Range("J89").Value = 401816
$ws.Range("K89").Value = 17395.3335
$ws.Range("L89").Value = 2009080
$ws.Range("M89").Value = -11779.3335
$ws.Range("N89").Value = -2020312
$ws.Range("H94").Value = 2045.75
$ws.Range("I94").Value = 592.5
$ws.Range("J94").Value = 3499
$ws.Range("K94").Value = 592.5
$ws.Range("L94").Value = 3499
$ws.Range("M94").Value = -141.5
$ws.Range("N94").Value = -4401
$ws.Range("H134").Value = 1535.0605
$ws.Range("I134").Value = 1226.7037
$ws.Range("J134").Value = 2922.6667
$ws.Range("K134").Value = 3680.1111
$ws.Range("L134").Value = 8768.000100000001
$ws.Range("M134").Value = -1145.1111
$ws.Range("N134").Value = -13838.0001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 466.05264
$ws.Range("I16").Value = 448.88235
$ws.Range("J16").Value = 612
$ws.Range("K16").Value = 448.88235
$ws.Range("L16").Value = 612
$ws.Range("M16").Value = -161.88235
$ws.Range("N16").Value = -1186
$ws.Range("H31").Value = 2601.3
$ws.Range("I31").Value = 1736.8572
$ws.Range("J31").Value = 3066.7693
$ws.Range("K31").Value = 1736.8572
$ws.Range("L31").Value = 3066.7693
$ws.Range("M31").Value = -1441.8572
$ws.Range("N31").Value = -3656.7693
$ws.Range("H34").Value = 2601.3
$ws.Range("I34").Value = 1736.8572
$ws.Range("J34").Value = 3066.7693
$ws.Range("K34").Value = 1736.8572
$ws.Range("L34").Value = 3066.7693
$ws.Range("M34").Value = -1534.8572
$ws.Range("N34").Value = -3470.7693
$ws.Range("H62").Value = 3500
$ws.Range("I62").Value = 3500
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3500
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2876
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 3500
$ws.Range("I65").Value = 3500
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 17500
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -14380
$ws.Range("N65").ClearContents()
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").ClearContents()
$ws.Range("N70").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").ClearContents()
$ws.Range("N73").Value = 0
$ws.Range("H113").Value = 466.05264
$ws.Range("I113").Value = 448.88235
$ws.Range("J113").Value = 612
$ws.Range("K113").Value = 448.88235
$ws.Range("L113").Value = 612
$ws.Range("M113").Value = 1721.11765
$ws.Range("N113").Value = -4952
$ws.Range("H141").Value = 64193.4
$ws.Range("I141").Value = 72000
$ws.Range("J141").Value = 62241.75
$ws.Range("K141").Value = 72000
$ws.Range("L141").Value = 62241.75
$ws.Range("M141").Value = -66820
$ws.Range("N141").Value = -72601.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 161.33333
$ws.Range("I12").Value = 47.666668
$ws.Range("J12").Value = 275
$ws.Range("K12").Value = 143.000004
$ws.Range("L12").Value = 825
$ws.Range("M12").Value = 29.99999600000001
$ws.Range("N12").Value = -1171
$ws.Range("H34").Value = 8961
$ws.Range("I34").Value = 14415.714
$ws.Range("J34").Value = 1324.4
$ws.Range("K34").Value = 43247.142
$ws.Range("L34").Value = 3973.2
$ws.Range("M34").Value = -43163.142
$ws.Range("N34").Value = -4141.200000000001
$ws.Range("H62").Value = 4000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 12000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -13372
$ws.Range("H65").Value = 4000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 36000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -42864
$ws.Range("H105").Value = 4771.2144
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 4771.2144
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 14313.6432
$ws.Range("N105").Value = -19555.6432
$ws.Range("H113").Value = 726.96155
$ws.Range("I113").Value = 841.75
$ws.Range("J113").Value = 706.0909
$ws.Range("K113").Value = 2525.25
$ws.Range("L113").Value = 2118.2727
$ws.Range("M113").Value = -355.25
$ws.Range("N113").Value = -6458.2727
$ws.Range("H131").Value = 14775.14
$ws.Range("I131").Value = 707.5
$ws.Range("J131").Value = 15998.413
$ws.Range("K131").Value = 2122.5
$ws.Range("L131").Value = 47995.239
$ws.Range("M131").Value = 2917.5
$ws.Range("N131").Value = -58075.239

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 30740.334
$ws.Range("I15").Value = 22222
$ws.Range("J15").Value = 34999.5
$ws.Range("K15").Value = 22222
$ws.Range("L15").Value = 34999.5
$ws.Range("M15").Value = -21934
$ws.Range("N15").Value = -35575.5
$ws.Range("H80").Value = 2420.6667
$ws.Range("I80").Value = 2738.5
$ws.Range("J80").Value = 1785
$ws.Range("K80").Value = 2738.5
$ws.Range("L80").Value = 1785
$ws.Range("M80").Value = -1740.5
$ws.Range("N80").Value = -3781
$ws.Range("H81").Value = 30740.334
$ws.Range("I81").Value = 22222
$ws.Range("J81").Value = 34999.5
$ws.Range("K81").Value = 22222
$ws.Range("L81").Value = 34999.5
$ws.Range("M81").Value = -21224
$ws.Range("N81").Value = -36995.5
$ws.Range("H83").Value = 2420.6667
$ws.Range("I83").Value = 2738.5
$ws.Range("J83").Value = 1785
$ws.Range("K83").Value = 13692.5
$ws.Range("L83").Value = 8925
$ws.Range("M83").Value = -8700.5
$ws.Range("N83").Value = -18909
$ws.Range("H84").Value = 30740.334
$ws.Range("I84").Value = 22222
$ws.Range("J84").Value = 34999.5
$ws.Range("K84").Value = 66666
$ws.Range("L84").Value = 104998.5
$ws.Range("M84").Value = -61674
$ws.Range("N84").Value = -114982.5
$ws.Range("H97").Value = 676.11536
$ws.Range("I97").Value = 685.9545000000001
$ws.Range("J97").Value = 622
$ws.Range("K97").Value = 685.9545000000001
$ws.Range("L97").Value = 622
$ws.Range("M97").Value = -189.9545000000001
$ws.Range("N97").Value = -1614
$ws.Range("H102").Value = 1727.5
$ws.Range("I102").Value = 1360.3334
$ws.Range("J102").Value = 2829
$ws.Range("K102").Value = 1360.3334
$ws.Range("L102").Value = 2829
$ws.Range("M102").Value = 261.6666
$ws.Range("N102").Value = -6073
$ws.Range("H113").Value = 712.1875
$ws.Range("I113").Value = 314
$ws.Range("J113").Value = 1224.1428
$ws.Range("K113").Value = 314
$ws.Range("L113").Value = 1224.1428
$ws.Range("M113").Value = 1856
$ws.Range("N113").Value = -5564.1428
$ws.Range("H126").Value = 2780691.2
$ws.Range("I126").Value = 3971274
$ws.Range("J126").Value = 2664.8333
$ws.Range("K126").Value = 11913822
$ws.Range("L126").Value = 7994.499899999999
$ws.Range("M126").Value = -11911352
$ws.Range("N126").Value = -12934.4999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2000
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2000
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1251
$ws.Range("H71").Value = 2000
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 10000
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -6256
$ws.Range("H132").Value = 2168.3125
$ws.Range("I132").Value = 1612.1875
$ws.Range("J132").Value = 2724.4375
$ws.Range("K132").Value = 4836.5625
$ws.Range("L132").Value = 8173.3125
$ws.Range("M132").Value = -2306.5625
$ws.Range("N132").Value = -13233.3125

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1871.8292
$ws.Range("I132").Value = 1183.5312
$ws.Range("J132").Value = 4319.1113
$ws.Range("K132").Value = 3550.5936
$ws.Range("L132").Value = 12957.3339
$ws.Range("M132").Value = -1020.5936
$ws.Range("N132").Value = -18017.3339
